# "Generate Report for Archive"
#
# The localization status report is regenerated: every cell whose status
# was "Ready for handoff" is now "In Translation", and the (now narrower)
# Status/zh-cn/de-de columns are re-sized to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# 1. Update status text wherever it said "Ready for handoff".
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# 2. Re-fit the columns that held that text to their new (narrower) width.
#    ColumnWidth is quantized by Excel to whole-pixel steps, so 12.5 is the
#    character-width value that lands on the same pixel width as the
#    regenerated report's column width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
